# adding averages and more checks
#
# Training Dashboard: refresh "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I)
# columns for rows 3-8 (new run date 16-Sep-2025, one week later -> periods
# shrink by 8 days each).
# Exam Dashboard: comments updated now that the exam date is valid, and the
# now-unused wide "COMMENTS" column is narrowed back down.
# Also: re-colour the title/header bands (white bold text) and drop the
# title row down to the normal font size.

$wb = $excel.ActiveWorkbook

$trainingWs = $wb.Worksheets.Item("Training Dashboard")
$examWs     = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------
# Training Dashboard: updated "PERIOD TO EXPIRE" / "LAST UPDATE" values
# ---------------------------------------------------------------------
$periodUpdates = @{
    3 = 334
    4 = 251
    5 = 210
    6 = 255
    7 = 254
    8 = -343
}

foreach ($row in $periodUpdates.Keys) {
    $trainingWs.Range("H$row").Value = $periodUpdates[$row]
}

# Write the new "LAST UPDATE" date as literal text (matching how the sheet
# already stores these dates) rather than letting it be auto-converted to
# a real date serial: build it as a text formula, then paste back as a
# value so the stored cell stays a plain string.
$lastUpdateRange = $trainingWs.Range("I3:I8")
$lastUpdateRange.Formula = '="16-Sep-2025"'
$lastUpdateRange.Copy() | Out-Null
$lastUpdateRange.PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Exam Dashboard: comments + column width
# ---------------------------------------------------------------------
$examWs.Range("E3").Value = "date is valid"
$examWs.Range("E4").Value = "date is valid"

# Column E ("COMMENTS") no longer needs to be extra wide now the remark is
# short; 14.2 characters of COM ColumnWidth rounds to the file's stored
# width of 15.
$examWs.Columns.Item(5).ColumnWidth = 14.2

# ---------------------------------------------------------------------
# Title (row 1) + header (row 2) restyle on both sheets: bold white text,
# title no longer oversized.
# ---------------------------------------------------------------------
foreach ($ws in @($trainingWs, $examWs)) {
    $ws.Range("A1").Font.Size = 11
    $ws.Range("A1").Font.Color = 16777215
    $ws.Rows.Item(2).Font.Color = 16777215
}
